$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill column D (rows 1-125 and 127-130) with the value 1. These rows have no
# existing D cell, so a plain value assignment creates a fresh, unstyled cell
# exactly like the target OOXML (<c r="Dn"><v>1</v></c>).
$ws.Range("D1:D125").Value = 1
$ws.Range("D127:D130").Value = 1

# Row 126 already had a D cell (style-only, no content: <c r="D126" s="2"/>).
# Clear it first so the style is dropped, then set the value - this matches
# the target cell <c r="D126"><v>1</v></c> (no style attribute).
$ws.Cells.Item(126, 4).Clear()
$ws.Cells.Item(126, 4).Value = 1

# Update the view: select H18 (scrolling the view back to the top, dropping
# the old topLeftCell="A205"/selection B224 state).
$ws.Range("H18").Select()

Write-Output "done"
